$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D column (Price) updates: use apostrophe-prefix to force text, then reset style ---
$ws.Range("D2").Value = "'55.889.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'2.512.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Value = "'490.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'140.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.513"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'2.511.16"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.0990"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Value = "'0.332"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Value = "'2.955.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'55.867.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'20.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.0000138"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'2.519.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'4.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'322.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'10.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'5.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'58.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Value = "'0.413"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Value = "'2.622.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "'7.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Value = "'0.0₃0799"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'150.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'18.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'1.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'5.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'3.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("D39").Value = "'34.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'0.615"
$ws.Range("D40").Style = "Normal"
$ws.Range("D43").Value = "'3.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'1.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'4.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'2.008.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'259.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.0913"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").Value = "'17.66"
$ws.Range("D51").Style = "Normal"

# --- E column (Volume/1h) updates ---
$ws.Range("E2").Value = "  +6.61%  "
$ws.Range("E3").Value = "  +8.31%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("E5").Value = "  +11.79%  "
$ws.Range("E6").Value = "  +16.22%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +9.75%  "
$ws.Range("E9").Value = "  +8.15%  "
$ws.Range("E10").Value = "  +12.37%  "
$ws.Range("E11").Value = "  +6.32%  "
$ws.Range("E12").Value = "  +9.08%  "
$ws.Range("E13").Value = "  +2.38%  "
$ws.Range("E14").Value = "  +7.78%  "
$ws.Range("E15").Value = "  +6.57%  "
$ws.Range("E16").Value = "  +11.08%  "
$ws.Range("E17").Value = "  +16.63%  "
$ws.Range("E18").Value = "  +7.02%  "
$ws.Range("E19").Value = "  +11.66%  "
$ws.Range("E20").Value = "  +7.98%  "
$ws.Range("E21").Value = "  +12.45%  "
$ws.Range("E23").Value = "  +12.72%  "
$ws.Range("E24").Value = "  +8.81%  "
$ws.Range("E25").Value = "  +14.51%  "
$ws.Range("E26").Value = "  +13.10%  "
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("E28").Value = "  +7.29%  "
$ws.Range("E29").Value = "  +8.24%  "
$ws.Range("E30").Value = "  +18.02%  "
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("E32").Value = "  +5.04%  "
$ws.Range("E33").Value = "  +8.02%  "
$ws.Range("E34").Value = "  +12.71%  "
$ws.Range("E35").Value = "  +9.58%  "
$ws.Range("E36").Value = "  +8.16%  "
$ws.Range("E39").Value = "  +8.70%  "
$ws.Range("E40").Value = "  +15.74%  "
$ws.Range("E43").Value = "  +9.84%  "
$ws.Range("E44").Value = "  +10.01%  "
$ws.Range("E45").Value = "  +14.37%  "
$ws.Range("E46").Value = "  +4.90%  "
$ws.Range("E47").Value = "  +35.37%  "
$ws.Range("E48").Value = "  +10.97%  "
$ws.Range("E51").Value = "  +13.29%  "

# --- Rows with reordered coins (B, C, D, E all updated) ---
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "'0.881"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.71%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'1.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.53%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "'0.0557"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.16%  "
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").Value = "'10.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0226"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.25%  "
